$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that sits right after the word
#    "Emotion" in the title line ("Project Title: Speech Based Summarization
#    and Emotion Analysis").
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Locate the paragraph that reads just "Evaluation: " near the end of the
#    document (section "6. Division of Labor between teammates"). There is
#    another, unrelated, "Evaluation: " paragraph earlier in the doc, so we
#    disambiguate using the preceding paragraph's text.
# ---------------------------------------------------------------------------
$target = $null
$paras = $d.Paragraphs
for ($i = 2; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.Trim() -eq "Evaluation:") {
        $prev = $paras.Item($i - 1)
        if ($prev.Range.Text.Trim().StartsWith("Summary generation for Evaluation")) {
            $target = $p
        }
    }
}

if ($target -ne $null) {
    $r = $target.Range

    # Find the ": " run inside this paragraph (colon followed by one space).
    $full = $d.Content.Text
    $colonIdx = $full.IndexOf(": ", $r.Start)

    $subRange = $d.Range($colonIdx, $colonIdx + 2)

    $rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="0E0E0E"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>'

    # Replace ": " with five runs: ":" / " Suraj, Raksha" / a "_GoBack"
    # bookmark / " " / " " -- matching the target markup exactly.
    $xmlSnippet = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r>$rPr<w:t>:</w:t></w:r>
<w:r>$rPr<w:t xml:space="preserve"> Suraj, Raksha</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r>$rPr<w:t xml:space="preserve"> </w:t></w:r>
<w:r>$rPr<w:t xml:space="preserve"> </w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

    $subRange.InsertXML($xmlSnippet)

    Write-Host "Updated paragraph text: [$($target.Range.Text)]"
}
else {
    Write-Host "WARNING: target 'Evaluation: ' paragraph not found"
}
